$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "65.336.88"
Set-TextValue $ws.Range("E2") "  -0.86%  "

Set-TextValue $ws.Range("D3") "3.329.47"
Set-TextValue $ws.Range("E3") "  -4.05%  "

Set-TextValue $ws.Range("E4") "  +0.00%  "

Set-TextValue $ws.Range("D5") "574.21"
Set-TextValue $ws.Range("E5") "  -1.27%  "

Set-TextValue $ws.Range("D6") "177.71"
Set-TextValue $ws.Range("E6") "  +2.95%  "

Set-TextValue $ws.Range("D7") "0.611"
Set-TextValue $ws.Range("E7") "  +1.98%  "

Set-TextValue $ws.Range("E8") "  +0.03%  "

Set-TextValue $ws.Range("D9") "3.328.32"
Set-TextValue $ws.Range("E9") "  -4.04%  "

Set-TextValue $ws.Range("E10") "  -1.22%  "

Set-TextValue $ws.Range("D11") "6.85"
Set-TextValue $ws.Range("E11") "  -0.04%  "

Set-TextValue $ws.Range("D12") "0.406"
Set-TextValue $ws.Range("E12") "  -0.81%  "

Set-TextValue $ws.Range("D13") "3.909.25"
Set-TextValue $ws.Range("E13") "  -4.04%  "

Set-TextValue $ws.Range("D14") "0.134"
Set-TextValue $ws.Range("E14") "  +0.39%  "

Set-TextValue $ws.Range("D15") "28.59"
Set-TextValue $ws.Range("E15") "  -4.32%  "

Set-TextValue $ws.Range("D16") "65.403.16"
Set-TextValue $ws.Range("E16") "  -0.92%  "

Set-TextValue $ws.Range("D17") "0.0000168"
Set-TextValue $ws.Range("E17") "  -1.22%  "

Set-TextValue $ws.Range("D18") "3.316.57"
Set-TextValue $ws.Range("E18") "  -4.49%  "

Set-TextValue $ws.Range("D19") "5.72"
Set-TextValue $ws.Range("E19") "  -3.19%  "

Set-TextValue $ws.Range("D20") "13.34"
Set-TextValue $ws.Range("E20") "  -3.70%  "

Set-TextValue $ws.Range("D21") "361.91"
Set-TextValue $ws.Range("E21") "  -1.06%  "

Set-TextValue $ws.Range("D22") "7.41"
Set-TextValue $ws.Range("E22") "  -3.81%  "

Set-TextValue $ws.Range("E23") "  -0.08%  "

Set-TextValue $ws.Range("D24") "71.26"
Set-TextValue $ws.Range("E24") "  -1.23%  "

Set-TextValue $ws.Range("D25") "0.518"
Set-TextValue $ws.Range("E25") "  -2.87%  "

Set-TextValue $ws.Range("D26") "0.0000121"
Set-TextValue $ws.Range("E26") "  -2.24%  "

Set-TextValue $ws.Range("D27") "9.49"
Set-TextValue $ws.Range("E27") "  -0.49%  "

Set-TextValue $ws.Range("E28") "  -0.90%  "

Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  +0.02%  "

Set-TextValue $ws.Range("E30") "  -1.03%  "

Set-TextValue $ws.Range("B31") "USDe"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D31") "0.999"
Set-TextValue $ws.Range("E31") "  -0.06%  "

Set-TextValue $ws.Range("B32") "NEARProtocol"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "5.59"
Set-TextValue $ws.Range("E32") "  -2.81%  "

Set-TextValue $ws.Range("D33") "22.85"
Set-TextValue $ws.Range("E33") "  -4.22%  "

Set-TextValue $ws.Range("D34") "6.81"
Set-TextValue $ws.Range("E34") "  -3.98%  "

Set-TextValue $ws.Range("E35") "  -5.65%  "

Set-TextValue $ws.Range("D36") "1.48"
Set-TextValue $ws.Range("E36") "  -2.95%  "

Set-TextValue $ws.Range("D37") "160.26"
Set-TextValue $ws.Range("E37") "  +0.68%  "

Set-TextValue $ws.Range("D38") "0.843"
Set-TextValue $ws.Range("E38") "  -4.97%  "

Set-TextValue $ws.Range("D39") "27.36"
Set-TextValue $ws.Range("E39") "  -6.38%  "

Set-TextValue $ws.Range("D40") "1.74"
Set-TextValue $ws.Range("E40") "  -0.45%  "

Set-TextValue $ws.Range("D41") "2.55"
Set-TextValue $ws.Range("E41") "  -0.15%  "

Set-TextValue $ws.Range("D42") "2.703.79"
Set-TextValue $ws.Range("E42") "  -3.59%  "

Set-TextValue $ws.Range("D43") "6.21"
Set-TextValue $ws.Range("E43") "  -4.12%  "

Set-TextValue $ws.Range("D44") "4.26"
Set-TextValue $ws.Range("E44") "  -3.67%  "

Set-TextValue $ws.Range("D45") "336.59"
Set-TextValue $ws.Range("E45") "  +7.27%  "

Set-TextValue $ws.Range("D46") "0.0666"
Set-TextValue $ws.Range("E46") "  -2.09%  "

Set-TextValue $ws.Range("E47") "  -1.00%  "

Set-TextValue $ws.Range("D48") "23.98"
Set-TextValue $ws.Range("E48") "  -0.36%  "

Set-TextValue $ws.Range("E49") "  -3.04%  "

Set-TextValue $ws.Range("E50") "  +1.84%  "

Set-TextValue $ws.Range("D51") "0.961"
Set-TextValue $ws.Range("E51") "  -1.19%  "
